$d = $word.ActiveDocument

# --- 1) Intro italic paragraph (Prueba 6 description) ---
$d.Content.Find.Execute(
    "El administrador no es capaz de editar la película que desea porque los datos colisionan con los de una película ya creada.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "El administrador desea eliminar una película de la base de datos.",
    2) | Out-Null

# --- 2) Criterio de aceptación bullet ---
$d.Content.Find.Execute(
    "La aplicación le debe preguntar al administrador si está seguro que desea guardar los cambios de la película ya que ya existe una con esos datos.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "La película no debe aparecer más en la lista de la cartelera.",
    2) | Out-Null

# --- 3) Pasos bullet 1 ---
$d.Content.Find.Execute(
    "El administrador llena los campos con los datos de la película. El nombre de la película es idéntico al de una ya creada.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "El administrador entra al menú de cartelera.",
    2) | Out-Null

# --- 4) Pasos bullet 2: "Luego hace click en el botón ..." ---
# This paragraph contains a middle run ("click") with italic formatting that
# must be preserved, so edit the surrounding runs via Range, not a single Find/Replace.
# There are several paragraphs with this same wording elsewhere in the document,
# so disambiguate by requiring the NEXT paragraph to be the "mensaje de
# confirmación" bullet that belongs to the same (Prueba 6) section.
for ($i = 1; $i -le $d.Paragraphs.Count - 1; $i++) {
    $para = $d.Paragraphs.Item($i)
    $nextPara = $d.Paragraphs.Item($i + 1)
    if ($para.Range.Text.StartsWith("Luego hace click") -and
        $nextPara.Range.Text.StartsWith("Aparece un mensaje de confirmación")) {
        $pRange = $para.Range
        $pStart = $pRange.Start
        $fullText = $pRange.Text
        $clickRelPos = $fullText.IndexOf("click")
        $clickStart = $pStart + $clickRelPos
        $clickEnd = $clickStart + 5

        # Text before "click": "Luego hace " -> "Luego navega por las páginas de la misma hasta elegir una película y hace "
        $beforeRange = $d.Range($pStart, $clickStart)
        $beforeRange.Text = "Luego navega por las páginas de la misma hasta elegir una película y hace "

        # Recompute paragraph/range after the edit above (length changed).
        $para2 = $d.Paragraphs.Item($i)
        $pRange2 = $para2.Range
        $fullText2 = $pRange2.Text
        $clickRelPos2 = $fullText2.IndexOf("click")
        $clickEnd2 = $pRange2.Start + $clickRelPos2 + 5

        # Text after "click" up to (not including) the paragraph mark.
        $afterRange = $d.Range($clickEnd2, $pRange2.End - 1)
        $afterRange.Text = " en el botón de eliminar que está al lado del botón de edición."

        break
    }
}

# --- 5) Pasos bullet 3 ---
$d.Content.Find.Execute(
    "Aparece un mensaje de confirmación y se le indica que la película ya existe.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Aparece un mensaje de confirmación.",
    2) | Out-Null

# --- 6) New Pasos bullet 4, appended after the "confirmación" bullet ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.StartsWith("Aparece un mensaje de confirmación.")) {
        $para.Range.InsertParagraphAfter()
        $newPara = $d.Paragraphs.Item($i + 1)
        $newPara.Range.Text = "Al afirmar la acción la película desaparece de la lista y de la base de datos."
        break
    }
}
